# Playtest Feedback.docx edit:
#   "buttons bigger" -> "buttons bigger (maybe? Made the entire pda bigger)"
#   "Code position indicator (while running)" -> "Code position indicator (while running) – Done (needs a test)"
#   "Forward until wall button" -> "Forward until wall button – Done (needs test)"
#
# Each addition is appended as new text right after the existing sentence,
# matching the commit "Button (and also pda) size increase".

$d = $word.ActiveDocument

# wdReplace = 2, wdFindContinue = 1
$wdReplaceAll = 2
$wdFindContinue = 1

$d.Content.Find.Execute(
    "buttons bigger",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "buttons bigger (maybe? Made the entire pda bigger)",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "Code position indicator (while running)",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Code position indicator (while running) – Done (needs a test)",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "Forward until wall button",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Forward until wall button – Done (needs test)",
    $wdReplaceAll)
